$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-order the country list: "Colombia" moves up so that it is listed
#    right after "Republica Dominicana" (it used to sit between "Panama" and
#    "Malasia"). "Australia" and "Panama" each shift down by one row as a
#    result, and each row now carries the figures for its (new) country.
# ---------------------------------------------------------------------------

# Row 47 (Republica Dominicana) is unchanged.

# Row 48 used to be Australia -> now Colombia, with Colombia's new figures.
$ws.Range("A48").Value = "Colombia"
$ws.Range("B48").Value = 7006
$ws.Range("C48").Value = 499
$ws.Range("D48").Value = 1551
$ws.Range("E48").Value = 5141
$ws.Range("F48").Value = 118
$ws.Range("G48").Value = 21
$ws.Range("H48").Value = 314

# Row 49 used to be Panama -> now Australia, carrying what used to be row
# 48's (Australia's) figures.
$ws.Range("A49").Value = "Australia"
$ws.Range("B49").Value = 6767
$ws.Range("C49").Value = 13
$ws.Range("D49").Value = 5745
$ws.Range("E49").Value = 929
$ws.Range("F49").Value = 28
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 93

# Row 50 used to be Colombia -> now Panama, carrying what used to be row
# 49's (Panama's) figures.
$ws.Range("A50").Value = "Panama"
$ws.Range("B50").Value = 6532
$ws.Range("C50").Value = 154
$ws.Range("D50").Value = 576
$ws.Range("E50").Value = 5768
$ws.Range("F50").Value = 86
$ws.Range("G50").Value = 10
$ws.Range("H50").Value = 188

# Row 51 (Malasia) is unchanged.

# ---------------------------------------------------------------------------
# 2. Bump the "last updated" timestamp banner in A1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 23:52"

# ---------------------------------------------------------------------------
# 3. Refresh the daily figures for the countries whose row did not move.
# ---------------------------------------------------------------------------

# Row 4 -> Estados Unidos
$ws.Range("B4").Value = 1125305
$ws.Range("C4").Value = 30282
$ws.Range("D4").Value = 160173
$ws.Range("E4").Value = 899592
$ws.Range("F4").Value = 16478
$ws.Range("G4").Value = 1684
$ws.Range("H4").Value = 65540

# Row 9 -> Alemania
$ws.Range("B9").Value = 164077
$ws.Range("C9").Value = 1068
$ws.Range("E9").Value = 30441
$ws.Range("G9").Value = 113
$ws.Range("H9").Value = 6736

# Row 61 -> Kazajistan
$ws.Range("B61").Value = 3597
$ws.Range("C61").Value = 195
$ws.Range("D61").Value = 922
$ws.Range("E61").Value = 2650

# Row 124 -> Venezuela
$ws.Range("B124").Value = 335
$ws.Range("C124").Value = 2
$ws.Range("E124").Value = 183
